$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2108.2307
$ws.Range("I98").Value = 2200.5833
$ws.Range("K98").Value = 2200.5833
$ws.Range("M98").Value = -702.5832999999998

$ws.Range("H101").Value = 2099.5
$ws.Range("I101").Value = 399
$ws.Range("K101").Value = 1197
$ws.Range("M101").Value = 425

$ws.Range("H122").Value = 2108.2307
$ws.Range("I122").Value = 2200.5833
$ws.Range("K122").Value = 6601.749899999999
$ws.Range("M122").Value = -4151.749899999999

$ws.Range("H138").Value = 3171.3572
$ws.Range("I138").Value = 3155.9
$ws.Range("J138").Value = 3210
$ws.Range("K138").Value = 9467.700000000001
$ws.Range("L138").Value = 9630
$ws.Range("M138").Value = -4327.700000000001
$ws.Range("N138").Value = -19910

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 7777
$ws.Range("J36").Value = 7777
$ws.Range("L36").Value = 7777
$ws.Range("N36").Value = -8469

$ws.Range("H102").Value = 1976.0834
$ws.Range("I102").Value = 1512.6666
$ws.Range("J102").Value = 3366.3333
$ws.Range("K102").Value = 1512.6666
$ws.Range("L102").Value = 3366.3333
$ws.Range("M102").Value = 109.3334
$ws.Range("N102").Value = -6610.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1766.6666
$ws.Range("I5").Value = 750
$ws.Range("J5").Value = 3800
$ws.Range("K5").Value = 750
$ws.Range("L5").Value = 3800
$ws.Range("M5").Value = -637
$ws.Range("N5").Value = -4026

$ws.Range("H105").Value = 1798.7084
$ws.Range("I105").Value = 1798.7084
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1798.7084
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -51.70839999999998
$ws.Range("N105").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 799.5
$ws.Range("I22").Value = 324.25
$ws.Range("J22").Value = 1274.75
$ws.Range("K22").Value = 324.25
$ws.Range("L22").Value = 1274.75
$ws.Range("M22").Value = 25.75
$ws.Range("N22").Value = -1974.75

$ws.Range("H58").Value = 1673650
$ws.Range("I58").Value = 3345625.8
$ws.Range("K58").Value = 3345625.8
$ws.Range("M58").Value = -3345422.8

$ws.Range("H86").Value = 58825300
$ws.Range("I86").Value = 90910610
$ws.Range("K86").Value = 90910610
$ws.Range("M86").Value = -90909487

$ws.Range("H89").Value = 58825300
$ws.Range("I89").Value = 90910610
$ws.Range("K89").Value = 454553050
$ws.Range("M89").Value = -454547434

$ws.Range("H107").Value = 361.1875
$ws.Range("I107").Value = 358.6
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 358.6
$ws.Range("L107").Value = 400
$ws.Range("M107").Value = 1561.4
$ws.Range("N107").Value = -4240

$ws.Range("H132").Value = 2033
$ws.Range("I132").Value = 1180.2727
$ws.Range("J132").Value = 4378
$ws.Range("K132").Value = 3540.8181
$ws.Range("L132").Value = 13134
$ws.Range("M132").Value = -1010.8181
$ws.Range("N132").Value = -18194

$ws.Range("H136").Value = 1673650
$ws.Range("I136").Value = 3345625.8
$ws.Range("K136").Value = 10036877.4
$ws.Range("M136").Value = -10034327.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 13665
$ws.Range("J9").Value = 16198
$ws.Range("L9").Value = 48594
$ws.Range("N9").Value = -49042

$ws.Range("H10").Value = 76.333336
$ws.Range("I10").Value = 76.333336
$ws.Range("K10").Value = 229.000008
$ws.Range("M10").Value = -90.00000800000001

$ws.Range("H126").Value = 4665
$ws.Range("I126").Value = 3330
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 9990
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -5050
$ws.Range("N126").Value = -27880

$ws.Range("H131").Value = 9402.704
$ws.Range("J131").Value = 10240.108
$ws.Range("L131").Value = 30720.324
$ws.Range("N131").Value = -40800.324

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2459.087
$ws.Range("I102").Value = 2446.1428
$ws.Range("K102").Value = 2446.1428
$ws.Range("M102").Value = -824.1428000000001

$ws.Range("H107").Value = 759.6
$ws.Range("I107").Value = 699.5
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 699.5
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 1220.5
$ws.Range("N107").Value = -4840

$ws.Range("H122").Value = 3000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = ""
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 10000
$ws.Range("I22").Value = 10000
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 10000
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -9705
$ws.Range("N22").Value = ""

$ws.Range("H27").Value = 10000
$ws.Range("I27").Value = 10000
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 10000
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -9893
$ws.Range("N27").Value = ""

$ws.Range("H40").Value = 12333.167
$ws.Range("I40").Value = 1000
$ws.Range("J40").Value = 14599.8
$ws.Range("K40").Value = 1000
$ws.Range("L40").Value = 14599.8
$ws.Range("M40").Value = -864
$ws.Range("N40").Value = -14871.8

$ws.Range("H68").Value = 4999.5
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 4999.5
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 4999.5
$ws.Range("M68").Value = ""
$ws.Range("N68").Value = -6497.5

$ws.Range("H71").Value = 4999.5
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 4999.5
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 24997.5
$ws.Range("M71").Value = ""
$ws.Range("N71").Value = -32485.5

$ws.Range("H122").Value = 7087
$ws.Range("I122").Value = 1802
$ws.Range("J122").Value = 9201
$ws.Range("K122").Value = 5406
$ws.Range("L122").Value = 27603
$ws.Range("M122").Value = -2956
$ws.Range("N122").Value = -32503

$ws.Range("H132").Value = 2355.5833
$ws.Range("I132").Value = 1999.6666
$ws.Range("K132").Value = 5998.9998
$ws.Range("M132").Value = -3468.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 269653.84
$ws.Range("I122").Value = 470019.25
$ws.Range("K122").Value = 1410019.75
$ws.Range("M122").Value = -1407607.75
